$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -49134425089.22117
$ws.Range("C2").Value = 1.809670881189074
$ws.Range("D2").Value = 0.009120199989411049

$ws.Range("B3").Value = -49436435013.78778
$ws.Range("C3").Value = 12.0644725412605
$ws.Range("D3").Value = 0.009230700001353398

$ws.Range("B4").Value = -49614087910.59167
$ws.Range("C4").Value = 18.09670881189074
$ws.Range("D4").Value = 0.008834100008243695

$ws.Range("B5").Value = -52634187156.25777
$ws.Range("C5").Value = 120.644725412605
$ws.Range("D5").Value = 0.009050199994817376

$ws.Range("B6").Value = -56187245092.33553
$ws.Range("C6").Value = 241.2894508252099
$ws.Range("D6").Value = 0.009175700004561804

$ws.Range("B7").Value = -56443775875.32035
$ws.Range("C7").Value = 250
$ws.Range("D7").Value = 0.008675800010678358

$ws.Range("B8").Value = -56443775875.32035
$ws.Range("C8").Value = 250
$ws.Range("D8").Value = 0.00929040000482928

$ws.Range("B9").Value = -56443775875.32035
$ws.Range("C9").Value = 250
$ws.Range("D9").Value = 0.01006959998630919

$ws.Range("B10").Value = -56443775875.32035
$ws.Range("C10").Value = 250
$ws.Range("D10").Value = 0.009266600012779236

$ws.Range("B11").Value = -56443775875.32035
$ws.Range("C11").Value = 250
$ws.Range("D11").Value = 0.009562000006553717

$ws.Range("D12").Value = 0.008792799999355339
